$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "#end" marker currently in A2, then rebuild the grid:
# A1:E2 all become "$" (same string used previously in B1), and the
# "#end" marker moves down to A3.
$ws.Range("A1:E2").Value = "$"
$ws.Range("A3").Value = "#end"

$ws.Range("A3").Select()
